$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header, formatted like the other header cells (bold,
# bordered, centered) by copying the format from the existing "sum" header.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Save values per row (0/1 flags), rows 2-11
$values = @(0, 0, 0, 0, 1, 0, 1, 0, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
